$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 21, shifting existing rows 21-26 down to 22-27
$ws.Range("A21:R21").Insert()

# Populate the new row 21 with the latest weekly data record
$ws.Range("A21").Value = 9
$ws.Range("B21").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C21").Value = "Metropolitana"
$ws.Range("D21").Value = 44460
$ws.Range("D21").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("E21").Value = 13
$ws.Range("F21").Value = 100112035
$ws.Range("G21").Value = "Bruselas (repollito)"
$ws.Range("H21").Value = "Sin especificar"
$ws.Range("I21").Value = "Primera"
$ws.Range("J21").Value = 25
$ws.Range("K21").Value = 24000
$ws.Range("L21").Value = 25000
$ws.Range("M21").Value = 24480
$ws.Range("N21").Value = "$/malla 15 kilos"
$ws.Range("O21").Value = "Hijuelas"
$ws.Range("P21").Value = 1632
$ws.Range("Q21").Value = 15
$ws.Range("R21").Value = "Hortaliza"
